# Append a new row (row 7) to the "سجل" sheet, recording a trip logged
# at 2025-05-01T11:17:19.059Z. Mirrors the existing rows: every cell is
# plain text (including numeric-looking values such as G7/H7), so numeric
# strings are entered with a leading apostrophe to force text storage and
# then restored to the default "Normal" style to avoid leaving a
# quote-prefix format behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "2025-05-01T11:17:19.059Z"
$ws.Range("B7").Value = "IDRF"
$ws.Range("C7").Value = "C3"
$ws.Range("D7").Value = "الرحلة 1"
$ws.Range("E7").Value = "الصمود"
$ws.Range("F7").Value = "يامن "

$ws.Range("G7").Value = "'123123"
$ws.Range("G7").Style = "Normal"

$ws.Range("H7").Value = "'123"
$ws.Range("H7").Style = "Normal"
